# Apply the "Trade #33 closed" update to the live trading results workbook.
#
# Changes:
#  1. Summary sheet           - refresh aggregate stats (capital, P&L, trade counts, win rate)
#  2. Strategy Status sheet   - refresh the MarketMaking strategy row
#  3. All Trades sheet        - append the new closed trade as row 34
#  4. MarketMaking sheet      - append the same closed trade as row 34

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.4    # Current Capital
$summary.Range("B4").Value = -0.6      # Total P&L $
$summary.Range("B5").Value = -0.36     # Total P&L %
$summary.Range("B6").Value = 33        # Total Trades
$summary.Range("B7").Value = 8         # Winning Trades
$summary.Range("B9").Value = 24.24     # Win Rate %

# ---------------------------------------------------------------------------
# 2. Strategy Status sheet (MarketMaking is row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.4       # Capital
$status.Range("D4").Value = 33         # Trades
$status.Range("E4").Value = -0.6       # P&L $
$status.Range("F4").Value = -0.6       # P&L %
$status.Range("G4").Value = 24.24      # Win Rate %

# ---------------------------------------------------------------------------
# Helper: write the newly-closed trade into row 34 of a trades-log sheet.
# Date/time columns must stay plain text, so force a text number format
# before assigning them (otherwise Excel auto-converts them to date/time
# serials).
# ---------------------------------------------------------------------------
function Add-Trade33Row([object]$ws) {
    $ws.Range("B34:C34").NumberFormat = "@"

    $ws.Range("A34").Value = 33
    $ws.Range("B34").Value = "2026-02-17"
    $ws.Range("C34").Value = "08:28:38"
    $ws.Range("D34").Value = "MarketMaking"
    $ws.Range("E34").Value = "UP"
    $ws.Range("F34").Value = 0.96
    $ws.Range("G34").Value = 0.99
    $ws.Range("H34").Value = "CLOSED"
    $ws.Range("I34").Value = 3.125
    $ws.Range("J34").Value = 0.03
    $ws.Range("K34").Value = 99.4
    $ws.Range("L34").Value = 0
    $ws.Range("M34").Value = 0
    $ws.Range("N34").Value = 0.6
    $ws.Range("O34").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P34").Value = "early_exit"
    $ws.Range("Q34").Value = 0.12
}

# ---------------------------------------------------------------------------
# 3. All Trades sheet
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-Trade33Row $allTrades

# ---------------------------------------------------------------------------
# 4. MarketMaking sheet
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-Trade33Row $marketMaking
